$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Missing CASE values (column A) discovered for rows 40-147.
$caseValues = @{
    40 = 231
    41 = 235
    42 = 236
    43 = 237
    44 = 241
    45 = 245
    46 = 246
    47 = 249
    48 = 250
    49 = 251
    50 = 253
    51 = 255
    52 = 257
    53 = 259
    54 = 262
    55 = 267
    56 = 269
    57 = 270
    58 = 271
    59 = 275
    60 = 278
    61 = 280
    62 = 284
    63 = 285
    64 = 287
    65 = 289
    66 = 290
    67 = 292
    68 = 293
    69 = 294
    70 = 299
    71 = 301
    72 = 302
    73 = 303
    74 = 305
    75 = 309
    76 = 318
    77 = 322
    78 = 323
    79 = 324
    80 = 325
    81 = 331
    82 = 332
    83 = 333
    84 = 335
    85 = 336
    86 = 338
    87 = 341
    88 = 343
    89 = 345
    90 = 346
    91 = 354
    92 = 357
    93 = 358
    94 = 359
    95 = 362
    96 = 365
    97 = 367
    98 = 371
    99 = 372
    100 = 375
    101 = 376
    102 = 378
    103 = 392
    104 = 402
    105 = 408
    106 = 409
    107 = 410
    108 = 411
    109 = 419
    110 = 420
    111 = 424
    112 = 425
    113 = 429
    114 = 437
    115 = 438
    116 = 439
    117 = 440
    118 = 448
    119 = 449
    120 = 450
    121 = 451
    122 = 452
    123 = 454
    124 = 456
    125 = 457
    126 = 458
    127 = 459
    128 = 460
    129 = 462
    130 = 463
    131 = 465
    132 = 467
    133 = 468
    134 = 470
    135 = 471
    136 = 472
    137 = 473
    138 = 474
    139 = 475
    140 = 476
    141 = 477
    142 = 478
    143 = 479
    144 = 480
    145 = 481
    146 = 491
    147 = 492

}

foreach ($row in $caseValues.Keys) {
    $val = $caseValues[$row]
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $val
    $cell.NumberFormat = "0"
}

# Reflect the reviewer's new scroll position / active selection.
$ws.Application.ActiveWindow.ScrollRow = 138
$ws.Range("K147").Select()
